# Applies the cryptos.xlsx crypto-price-table update described by the commit diff.
# Values are written via Range.Value; cells whose new content looks like a plain
# number (e.g. "226.21") are temporarily forced to Text format so Excel keeps them
# as strings (matching the original inlineStr cell type) and the Text format is
# cleared again right after so no stray cell style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.300.82"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.790.79"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.34"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "2.049.03"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.793.09"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "34.275.76"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.87"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.94%  "
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "1.441.24"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  +8.41%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.924"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "1.945.06"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  -5.89%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
